$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# This shared string is referenced from the Overview sheet (zh-cn / de-de
# status columns) as well as the per-locale "Status" column, so updating
# each referencing cell keeps them all in sync via the shared string table.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "Status"-related columns ---
# Target stored width is ~13.41 character-units; ColumnWidth snaps to the
# nearest whole pixel internally, so 12.5 is the closest achievable value.
$targetColumnWidth = 12.5

# Overview: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = $targetColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $targetColumnWidth

# zh-cn / de-de: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = $targetColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $targetColumnWidth
